$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the target paragraph: "v. Scopul senzorilor în cadrul ..."
# ------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.IndexOf("Scopul senzorilor") -ge 0) {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIndex)

# ------------------------------------------------------------------
# 1. Indentation: left=1440 -> left=1800, hanging=360
# ------------------------------------------------------------------
$p.Range.ParagraphFormat.LeftIndent = 90
$p.Range.ParagraphFormat.FirstLineIndent = -18

# ------------------------------------------------------------------
# 2. Split "v. Scopul senzorilor în cadrul aplicației: " into
#    "v. " + TAB + "Scopul senzorilor în cadrul aplicației: "
# ------------------------------------------------------------------
$oldHeading = "v. Scopul senzorilor în cadrul aplicației: "
$newHeading = "v. ^tScopul senzorilor în cadrul aplicației: "
$p.Range.Find.Execute($oldHeading, $true, $false, $false, $false, $false, $true, 1, $false, $newHeading, 2) | Out-Null

# ------------------------------------------------------------------
# 3. Append new sentence after "Utilizatorul va "
# ------------------------------------------------------------------
$oldTail = "Utilizatorul va "
$newTail = "Utilizatorul va avea datele afișate pe telefon, iar apăsând butonul de find-spot se va folosi geolocația lui pentru a găsi cel mai apropiat loc de parcare liber."
$p.Range.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null

# ------------------------------------------------------------------
# 4. Insert a new, empty paragraph right after this one
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item($targetIndex + 1)
$p2.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.ParagraphFormat.Alignment = 3
$newPara.Range.Font.Size = 12
$newPara.Range.Font.SizeBi = 12
